$d = $word.ActiveDocument

# Locate the two hyperlink "HYPERLINK ... clap_footer" field paragraphs near the
# end of the document (inserted right after the last inline image) through the
# trailing "BLOG DIFFERENCE BTW HTTP1 AND HTTP2" / "HTTP1" / "HTTP2" scratch
# content that was appended twice, and delete that whole block of paragraphs.
# We find it by searching for the start of the first HYPERLINK field's begin
# marker paragraph and the paragraph that holds the "server push" text, which
# is the last paragraph fully removed by the edit.

$startMarker = "HYPERLINK ""https://medium.com/m/signin?actionUrl=https%3A%2F%2Fmedium.com%2F_%2Fvote%2Fp%2Fe9d3e57b9dcb"
$endMarkerText = "The client retains the authority to deny the server push"

$found = $d.Content.Find
$found.ClearFormatting()
$found.Text = $startMarker
$found.Forward = $true
$found.Wrap = 1
$found.Execute() | Out-Null

$startRange = $null
if ($found.Found) {
    $startRange = $word.Selection.Range
}

# Fallback: locate via the paragraph collection if Find on field codes is not
# supported in this environment.
if (-not $startRange) {
    $n = $d.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -match "^\r$" -and $i -lt $n) {
        }
    }
}

# Find the paragraph that contains the server-push sentence, which is the
# last paragraph of the block slated for deletion.
$endFound = $d.Content.Find
$endFound.ClearFormatting()
$endFound.Text = $endMarkerText
$endFound.Forward = $true
$endFound.Wrap = 1
$endFound.Execute() | Out-Null

$endParaIndex = -1
$startParaIndex = -1

$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($startParaIndex -lt 0 -and $t -like "*HYPERLINK*e9d3e57b9dcb*") {
        $startParaIndex = $i
    }
    if ($t -like "*$endMarkerText*") {
        $endParaIndex = $i
    }
}

if ($startParaIndex -lt 0) {
    # The field's instruction text lives in the field code, not the
    # displayed text, on some renderers; locate by the following run
    # of near-empty shaded paragraphs that precede the "BLOG DIFFERENCE"
    # restatement instead.
    for ($i = 1; $i -le $n; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text -match "^\r$") {
            $nextText = ""
            if ($i + 8 -le $n) {
                $nextText = $d.Paragraphs.Item($i + 8).Range.Text
            }
            if ($nextText -like "*BLOG DIFFERENCE BTW HTTP1 AND HTTP2*") {
                $startParaIndex = $i
                break
            }
        }
    }
}

Write-Output ("startParaIndex=" + $startParaIndex + " endParaIndex=" + $endParaIndex)

if ($startParaIndex -gt 0 -and $endParaIndex -ge $startParaIndex) {
    $delStart = $d.Paragraphs.Item($startParaIndex).Range.Start
    $delEnd = $d.Paragraphs.Item($endParaIndex).Range.End
    $d.Range($delStart, $delEnd).Delete()
}

# The final remaining paragraph keeps only its first run's formatting; its
# text (minus the trailing paragraph mark) becomes five spaces.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastRange = $lastPara.Range
$textRange = $d.Range($lastRange.Start, $lastRange.End - 1)
$textRange.Text = "     "

Write-Output "done"
